$d = $word.ActiveDocument

# Locate the paragraph that ends with "python manage.py runserver"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd("`r", "`n") -eq "python manage.py runserver") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    # Insert a new empty paragraph right after it
    $target.Range.InsertParagraphAfter()

    # The paragraph following the one we just created is the new empty paragraph;
    # insert another paragraph after that one and fill it with "deactivate"
    $emptyPara = $target.Next()
    $emptyPara.Range.InsertParagraphAfter()
    $deactivatePara = $emptyPara.Next()
    $deactivatePara.Range.Text = "deactivate"
}
